$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 19
$ws.Range("D8").Value = "'2"
$ws.Range("E8").Value = "Short point (up to 3 mtr.)"
$ws.Range("F8").Value = 256
$ws.Range("G8").Value = "'4864.00"
$ws.Range("A9").Value = "P. point"
$ws.Range("C9").Value = 96
$ws.Range("D9").Value = "'4"
$ws.Range("E9").Value = "Long point  (up to 10 mtr.)"
$ws.Range("F9").Value = 662
$ws.Range("G9").Value = "'63552.00"
$ws.Range("A10").Value = "'"
$ws.Range("C10").Value = 75
$ws.Range("D10").Value = "'2.0"
$ws.Range("E10").Value = "Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it's  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = "'0.00"
$ws.Range("A11").Value = "Mtr."
$ws.Range("C11").Value = 71
$ws.Range("D11").Value = "'19"
$ws.Range("E11").Value = "2 x 2.5 sq. mm. + 1x1.5sqmm"
$ws.Range("F11").Value = 81
$ws.Range("G11").Value = "'5751.00"
$ws.Range("A12").Value = "Set"
$ws.Range("C12").Value = 48
$ws.Range("D12").Value = "'13.0"
$ws.Range("E12").Value = "Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. 'B' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure ""A"" attached with this BSR .   "
$ws.Range("F12").Value = 5733
$ws.Range("G12").Value = "'275184.00"
$ws.Range("A13").Value = "Each"
$ws.Range("C13").Value = 64
$ws.Range("D13").Value = "'25"
$ws.Range("E13").Value = "1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )"
$ws.Range("F13").Value = 1890
$ws.Range("G13").Value = "'120960.00"
$ws.Range("C14").Value = 12
$ws.Range("D14").Value = "'16.0"
$ws.Range("E14").Value = "Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure ""A"" attached with this BSR ."
$ws.Range("A15").Value = "Each"
$ws.Range("C15").Value = 87
$ws.Range("D15").Value = "'27"
$ws.Range("E15").Value = "1170mm(+/-10%) LED batten with min. lumen output 2200 lm"
$ws.Range("F15").Value = 492
$ws.Range("G15").Value = "'42804.00"
$ws.Range("C16").Value = 29
$ws.Range("A17").Value = "'"
$ws.Range("C17").Value = 96
$ws.Range("D17").Value = "'31"
$ws.Range("E17").Value = "Double pole MCB(With B/C curve tripping Characteristics)"
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = "'0.00"
$ws.Range("C18").Value = 89
$ws.Range("D18").Value = "'32"
$ws.Range("E18").Value = " 50/63 A rating"
$ws.Range("F18").Value = 900
$ws.Range("G18").Value = "'80100.00"
$ws.Range("A19").Value = "%"
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = "'37"
$ws.Range("E19").Value = "Add Tender Premium "
$ws.Range("G21").Value = "'593215.00"
$ws.Range("H21").Value = "'593215.00"
$ws.Range("G23").Value = "'593215.00"
$ws.Range("H23").Value = "'593215.00"
